# Update Leave Card 4/12/2023 4:43 PM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 64: SP(1-0-0) entry now has an EARNED value of 1.25 ---
$ws.Range("C64").Value = 1.25

# --- Row 65: new SL(1-0-0) entry ---
$ws.Range("A65").Value = 44986
$ws.Range("B65").Value = "SL(1-0-0)"
$ws.Range("C65").Value = 1.25
$ws.Range("H65").Value = 1

# K65 gets a remark date, formatted like the other date remarks in column K (copy format from K64)
$ws.Range("K64").Copy()
$ws.Range("K65").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("K65").Value = 44999

# --- Row 66: new SP(1-0-0) entry (remarks only) ---
$ws.Range("B66").Value = "SP(1-0-0)"
$ws.Range("K64").Copy()
$ws.Range("K66").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("K66").Value = 44988

# --- Rows 67-91: fill in the PERIOD (column A) dates for the following months ---
$periodDates = @{
    67 = 45017
    68 = 45047
    69 = 45078
    70 = 45108
    71 = 45139
    72 = 45170
    73 = 45200
    74 = 45231
    75 = 45261
    76 = 45292
    77 = 45323
    78 = 45352
    79 = 45383
    80 = 45413
    81 = 45444
    82 = 45474
    83 = 45505
    84 = 45536
    85 = 45566
    86 = 45597
    87 = 45627
    88 = 45658
    89 = 45689
    90 = 45717
    91 = 45748
}
foreach ($r in $periodDates.Keys) {
    $ws.Range("A$r").Value = $periodDates[$r]
}

# --- Table1: insert one more row at the bottom (grows from A8:K139 to A8:K140) ---
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# The brand-new physical row (140) should carry the calculated-column formula the
# same way every other table row does.
$ws.Range("G140").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# ListRows.Add() inserts the new row before the previous last row, but leaves the
# *formatting* of the old "last row" (139) in place and puts plain formatting on
# the new physical row (140). Swap the formatting back: row 139 should look like a
# normal data row (copy format from row 138) and row 140 should carry the special
# "last row" formatting that row 139 used to have.
$ws.Range("A139:K139").Copy()
$ws.Range("A140:K140").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A138:K138").Copy()
$ws.Range("A139:K139").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- View state: the active selection moved down to B67, and the split moved up one row ---
$ws.Activate()
$excel.ActiveWindow.SplitRow = 54
$ws.Range("B67").Select()
